# Remove the "SELECTED WORK / PORTFOLIO" section (heading + its three
# bullet items) that follows the "Minor: Mathematics, Philosophy" line
# at the end of the EDUCATION section.

$d = $word.ActiveDocument

$headingIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "SELECTED WORK / PORTFOLIO") {
        $headingIndex = $i
    }
    $i = $i + 1
}

if ($headingIndex -ge 1) {
    $startPara = $d.Paragraphs($headingIndex)
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)

    $start = $startPara.Range.Start
    $end = $lastPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
